$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was added to the daily logic subset.
# It is inserted as row 131, pushing the existing rows 131-191 down to 132-192.
$ws.Rows(131).Insert()

# Populate the newly inserted row 131 with the new record.
# (Same underlying data as the old row 131 except for the date and volume,
# matching how the source daily feed appends/re-sorts records.)
$ws.Cells.Item(131, 1).Value2 = 10
$ws.Cells.Item(131, 2).Value2 = "Vega Modelo de Temuco"
$ws.Cells.Item(131, 3).Value2 = "La Araucanía"
$ws.Cells.Item(131, 4).Value2 = 44489
$ws.Cells.Item(131, 5).Value2 = 9
$ws.Cells.Item(131, 6).Value2 = 100112017
$ws.Cells.Item(131, 7).Value2 = "Apio"
$ws.Cells.Item(131, 8).Value2 = "Americana (o)"
$ws.Cells.Item(131, 9).Value2 = "Primera"
$ws.Cells.Item(131, 10).Value2 = 65
$ws.Cells.Item(131, 11).Value2 = 8000
$ws.Cells.Item(131, 12).Value2 = 8000
$ws.Cells.Item(131, 13).Value2 = 8000
$ws.Cells.Item(131, 14).Value2 = "`$/docena de matas"
$ws.Cells.Item(131, 15).Value2 = "Provincia del Elquí"
$ws.Cells.Item(131, 16).Value2 = 1333
$ws.Cells.Item(131, 17).Value2 = 6
$ws.Cells.Item(131, 18).Value2 = "Hortaliza"
